# Applies updated NATMI LR-pair values ("Natmi following Dr Hou advice")
# Columns order: row, E, G, H, I, J, K, M, N, O, P, Q, R, S, T
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 138.2190853333333, 414.657256, 0.2249223651785973, 0.2476599003709697, 13.63614208264534, 122.725278743808, 0.002270152386121829, 0.003686772552152797),
    @(3, 3, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 147.91433, 443.74299, 0.2406993279341593, 0.2650317656414439, 14.59263614048, 131.33372526432, 0.002429390039597752, 0.003945377663286839),
    @(4, 3, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 74.27261733333333, 222.817852, 0.1208629960061633, 0.1330811078998542, 7.327439335637333, 65.946954020736, 0.001219876105070113, 0.001981103017001694),
    @(5, 3, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 3, 84.85695366666668, 254.570861, 0.138086767645209, 0.1520460408212704, 8.371647620938669, 75.34482858844801, 0.001393716471070842, 0.002263423223233564),
    @(6, 3, 0.09865600000000001, 0.295968, 0.01009304870291239, 0.01488643315542961, 2, 169.2560955, 338.512191, 0.275428543235871, 0.2021811852664618, 16.698129357648, 100.188776145888, 0.002779913701051858, 0.003009756699754712),
    @(7, 3, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 138.2190853333333, 414.657256, 0.2249223651785973, 0.2476599003709697, 32.31212060077156, 290.8090854069441, 0.005379339496312122, 0.008736154156415585),
    @(8, 3, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 147.91433, 443.74299, 0.2406993279341593, 0.2650317656414439, 34.57862319097333, 311.20760871876, 0.00575666808618112, 0.009348943278756417),
    @(9, 3, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 74.27261733333333, 222.817852, 0.1208629960061633, 0.1330811078998542, 17.36305635956089, 156.267507236048, 0.002890611111715428, 0.00469440984260358),
    @(10, 3, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 3, 84.85695366666668, 254.570861, 0.138086767645209, 0.1520460408212704, 19.83740605777378, 178.536654519964, 0.003302542201715345, 0.005363394112238674),
    @(11, 3, 0.2337746666666667, 0.7013240000000001, 0.02391642775070728, 0.03527480283780177, 2, 169.2560955, 338.512191, 0.275428543235871, 0.2021811852664618, 39.567787306814, 237.406723840884, 0.006587266854783265, 0.007131901447787511),
    @(12, 2, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 138.2190853333333, 414.657256, 0.2249223651785973, 0.2476599003709697, 1305.094666368393, 7830.56799821036, 0.2172728732961634, 0.2352369736624013),
    @(13, 2, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 147.91433, 443.74299, 0.2406993279341593, 0.2650317656414439, 1396.639275226775, 8379.835651360649, 0.2325132698083804, 0.2517374446994006),
    @(14, 2, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 74.27261733333333, 222.817852, 0.1208629960061633, 0.1330811078998542, 701.2982071556033, 4207.78924293362, 0.1167525087893778, 0.1264055950402489),
    @(15, 2, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 3, 84.85695366666668, 254.570861, 0.138086767645209, 0.1520460408212704, 801.2378129080893, 4807.426877448535, 0.1333905089724228, 0.1444192234857982),
    @(16, 2, 9.4422175, 18.884435, 0.9659905235463803, 0.9498387640067686, 2, 169.2560955, 338.512191, 0.275428543235871, 0.2021811852664618, 1598.152866911771, 6392.611467647085, 0.2660613626800359, 0.1920395271189196)
)

# Column letters for the values following the row number in each entry above
$cols = @("E", "G", "H", "I", "J", "K", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$colLetter$r").Value = $value
    }
}
